$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each updated cell as literal text, preserving the inline-string/text
# cell type (these columns store numbers-as-text, percents, names and URLs).
# NumberFormat '@' forces text interpretation so numeric-looking values like
# '300.66' or '1.00' are not silently coerced into floating point numbers;
# resetting the Style back to 'Normal' afterwards keeps the cell's style index
# identical to the original (unstyled) cell.
$cellUpdates = @(
    @{ Cell = 'D2'; Value = '45.696.04' }
    @{ Cell = 'E2'; Value = '  -2.03%  ' }
    @{ Cell = 'D3'; Value = '2.415.79' }
    @{ Cell = 'E3'; Value = '  +5.10%  ' }
    @{ Cell = 'E4'; Value = '  +0.02%  ' }
    @{ Cell = 'D5'; Value = '300.66' }
    @{ Cell = 'E5'; Value = '  -0.94%  ' }
    @{ Cell = 'D6'; Value = '97.08' }
    @{ Cell = 'E6'; Value = '  -3.76%  ' }
    @{ Cell = 'D7'; Value = '0.566' }
    @{ Cell = 'E7'; Value = '  -0.04%  ' }
    @{ Cell = 'E8'; Value = '  +0.04%  ' }
    @{ Cell = 'D9'; Value = '0.515' }
    @{ Cell = 'E9'; Value = '  -1.29%  ' }
    @{ Cell = 'D10'; Value = '34.98' }
    @{ Cell = 'E10'; Value = '  -4.56%  ' }
    @{ Cell = 'E11'; Value = '  +0.29%  ' }
    @{ Cell = 'D12'; Value = '7.18' }
    @{ Cell = 'E12'; Value = '  -2.61%  ' }
    @{ Cell = 'E13'; Value = '  +1.00%  ' }
    @{ Cell = 'D14'; Value = '2.785.80' }
    @{ Cell = 'E14'; Value = '  +5.25%  ' }
    @{ Cell = 'D15'; Value = '2.398.77' }
    @{ Cell = 'E15'; Value = '  +4.51%  ' }
    @{ Cell = 'D16'; Value = '14.34' }
    @{ Cell = 'E16'; Value = '  +3.75%  ' }
    @{ Cell = 'D17'; Value = '0.849' }
    @{ Cell = 'E17'; Value = '  +4.49%  ' }
    @{ Cell = 'D18'; Value = '45.678.86' }
    @{ Cell = 'E18'; Value = '  -1.98%  ' }
    @{ Cell = 'D19'; Value = '13.21' }
    @{ Cell = 'E19'; Value = '  +0.62%  ' }
    @{ Cell = 'D20'; Value = '0.0₃0953' }
    @{ Cell = 'E20'; Value = '  +1.66%  ' }
    @{ Cell = 'D21'; Value = '6.24' }
    @{ Cell = 'E21'; Value = '  +4.19%  ' }
    @{ Cell = 'D22'; Value = '67.46' }
    @{ Cell = 'E22'; Value = '  +1.53%  ' }
    @{ Cell = 'D23'; Value = '243.60' }
    @{ Cell = 'E23'; Value = '  -2.13%  ' }
    @{ Cell = 'E24'; Value = '  -2.60%  ' }
    @{ Cell = 'D25'; Value = '1.95' }
    @{ Cell = 'E25'; Value = '  +0.91%  ' }
    @{ Cell = 'D26'; Value = '1.00' }
    @{ Cell = 'E26'; Value = '  -0.08%  ' }
    @{ Cell = 'D27'; Value = '38.92' }
    @{ Cell = 'E27'; Value = '  -9.20%  ' }
    @{ Cell = 'E28'; Value = '  -1.95%  ' }
    @{ Cell = 'E29'; Value = '  +0.50%  ' }
    @{ Cell = 'B30'; Value = 'EthereumClassic' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' }
    @{ Cell = 'D30'; Value = '21.49' }
    @{ Cell = 'E30'; Value = '  +7.39%  ' }
    @{ Cell = 'B31'; Value = 'LidoDAOToken' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'D31'; Value = '3.84' }
    @{ Cell = 'E31'; Value = '  +17.20%  ' }
    @{ Cell = 'D32'; Value = '5.59' }
    @{ Cell = 'E32'; Value = '  -1.11%  ' }
    @{ Cell = 'E33'; Value = '  -1.62%  ' }
    @{ Cell = 'D34'; Value = '148.34' }
    @{ Cell = 'E34'; Value = '  +0.96%  ' }
    @{ Cell = 'D35'; Value = '0.0779' }
    @{ Cell = 'E35'; Value = '  -2.12%  ' }
    @{ Cell = 'D36'; Value = '2.00' }
    @{ Cell = 'E36'; Value = '  +12.33%  ' }
    @{ Cell = 'E37'; Value = '  -0.61%  ' }
    @{ Cell = 'E38'; Value = '  -1.41%  ' }
    @{ Cell = 'D39'; Value = '15.33' }
    @{ Cell = 'E39'; Value = '  -4.11%  ' }
    @{ Cell = 'D40'; Value = '3.92' }
    @{ Cell = 'E40'; Value = '  -1.86%  ' }
    @{ Cell = 'E41'; Value = '  -0.05%  ' }
    @{ Cell = 'D42'; Value = '3.28' }
    @{ Cell = 'E42'; Value = '  -2.14%  ' }
    @{ Cell = 'D43'; Value = '1.951.97' }
    @{ Cell = 'E43'; Value = '  +7.79%  ' }
    @{ Cell = 'E44'; Value = '  +0.10%  ' }
    @{ Cell = 'D45'; Value = '91.35' }
    @{ Cell = 'E45'; Value = '  +3.88%  ' }
    @{ Cell = 'E46'; Value = '  -10.38%  ' }
    @{ Cell = 'D47'; Value = '8.67' }
    @{ Cell = 'E47'; Value = '  +9.66%  ' }
    @{ Cell = 'D48'; Value = '15.41' }
    @{ Cell = 'E48'; Value = '  +14.30%  ' }
    @{ Cell = 'D49'; Value = '102.74' }
    @{ Cell = 'E49'; Value = '  +7.07%  ' }
    @{ Cell = 'D50'; Value = '0.188' }
    @{ Cell = 'E50'; Value = '  -3.19%  ' }
    @{ Cell = 'D51'; Value = '2.656.43' }
    @{ Cell = 'E51'; Value = '  +5.27%  ' }
)

foreach ($update in $cellUpdates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.Style = "Normal"
}
